# "add badge frontend to the the github remie"
#
# Appends 14 badge/ID-card payload rows (Code-39 style fixed-width strings,
# one per scanned badge) below the existing movie table in column A,
# widens the columns so the long badge text and the existing data stay
# readable, sets the print setup the author ended up with, and leaves the
# selection where it was when the workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badgeValues = @(
    '2500000100048004NON10O003400000NN0O 01320000000020200100                        1333100          '
    '2500000200048004NON10O003400000NN0O 01320002000020200100                        1333100          '
    '2500000300048004NON10O003400000NN0O 01320002222220200100                        1333100          '
    '2500000400048004NON10O003400000NN0O 01320002222220200100                        1333100          '
    '2500000500048004NON10O003400000NN0O 01320002202220000100                        1333100          '
    '2500000600048004NON10O003400000NN0O 01320002202220000100                        1333100          '
    '2500000700048004NON10O003400000NN0O 01320002222222200100                        1333100          '
    '2500000800048004NON10O003400000NN0O 01320002202220000100                        1333100          '
    '2500000900048004NON10O003400000NN0O 01320002202220000100                        1333100          '
    '2500001000048004NON10O003400000NN0O 01320002202220000100                        1333100          '
    '2500002000048004NON10O003400000NN0O 01320002202222200100                        1333100          '
    '2500004000048004NON10O003400000NN0O 01320002202220000100                        1333100          '
    '2500004100048004NON10O003400000NN0O 01320002202220000100                        1333100          '
    '2500000100049125NON10O003400000NN0O 01320002202220000100                        1333100          '
)

$startRow = 10
for ($i = 0; $i -lt $badgeValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $badgeValues[$i]
}

# Column A needs to be a lot wider to show the full badge payload; B and D
# get nudged slightly wider too, matching the resized widths in the file.
$ws.Columns.Item(1).ColumnWidth = 84
$ws.Columns.Item(2).ColumnWidth = 9
$ws.Columns.Item(4).ColumnWidth = 5.7

# Print setup the workbook was saved with: A4, portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it before saving.
$ws.Range("F15").Select() | Out-Null
